$d = $word.ActiveDocument

# 1. Update the Ativacao date from 2018 to 2022
$d.Content.Find.Execute("Ativação: 01/01/2018", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "Ativação: 01/01/2022", 2) | Out-Null

# Helper function: find the paragraph whose text contains $marker and insert a new
# italic paragraph right after it containing $newText.
function Insert-ItalicParagraphAfter($marker, $newText) {
    $paras = $d.Paragraphs
    $count = $paras.Count
    for ($i = 1; $i -le $count; $i++) {
        $p = $paras.Item($i)
        $t = $p.Range.Text
        if ($t -like "*$marker*") {
            $p.Range.InsertParagraphAfter()
            $newPara = $d.Paragraphs.Item($i + 1)
            $r = $newPara.Range
            $r.Text = $newText
            $r2 = $newPara.Range
            $r2.End = $r2.End - 1
            $r2.Font.Italic = $true
            return $true
        }
    }
    return $false
}

# 2. Insert English translation paragraph after the "Objetivos" paragraph text
Insert-ItalicParagraphAfter "Proporcionar aos alunos conhecimentos sobre a gestão ambiental" `
    "Provide students with knowledge about environmental management in companies, environmental policies, environmental management systems (EMS) and ISO 14000 series standards, enabling them to participate in the planning and implementation of an EMS in a company." | Out-Null

# 3. Insert English translation paragraph after the "Programa resumido" paragraph text
Insert-ItalicParagraphAfter "Sistemas de Gestão Ambiental; Iso 14000; Auditoria Ambiental." `
    "Environmental Management Systems; Iso 14000; Environmental Audit." | Out-Null

# 4. Replace the "Programa" paragraph text with the new expanded content
$oldPrograma = "Evolução das práticas de gestão ambiental empresarial; Normas (série ISO 14.000) e certificações ambientais; Produção mais limpa; Ferramentas de gestão ambiental focadas no produto: análise do ciclo de vida, ecodesign e rotulagem ambiental; Implantação do sistema de gerenciamento ambiental (SGA); inovação e sustentabilidade; otimização do ciclo de vida"
$newPrograma = "Evolução das práticas de gestão ambiental empresarial;- Economia circular, conceitos e aplicações;- Responsabilidade social corporativa: conceito e programa;- Implantação do sistema de gerenciamento ambiental (SGA): conceitos e modelos;- Produção mais limpa;- Ferramentas de gestão focadas no produto;- Análise e otimização do ciclo de vida do produto;- Ecoinovação e Ecodesign;- Rotulagem ambiental;- Inovação e sustentabilidade;- Normas ISO 14001 (série ISO 14000), requisitos e orientações para uso e Certificações ambientais."

$d.Content.Find.Execute($oldPrograma, $true, $false, $false, $false, $false, `
                         $true, 1, $false, $newPrograma, 2) | Out-Null

# 5. Insert English translation paragraph after the new "Programa" paragraph text
Insert-ItalicParagraphAfter "Normas ISO 14001 (série ISO 14000), requisitos e orientações para uso e Certificações ambientais." `
    "Evolution of corporate environmental management practices;- Circular economy, concepts and applications;- Corporate social responsibility: concept and program;- Implementation of the environmental management system (SGA): concepts and models;- Cleaner production;- Management tools focused on the product;- Analysis and optimization of the product life cycle;- Eco-innovation and Ecodesign;- Environmental labeling;- Innovation and sustainability;- ISO 14001 standards (ISO 14000 series), requirements and guidelines for use and Environmental Certifications." | Out-Null

Write-Host "Done"
